# "added 4wk low sales check"
# Updates forecast figures (MyForecast, Inventory Coverage, Seasonality Index)
# on the "Forecast Comparison" sheet, and the derived 4/8/16-week forecast
# totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Forecast Comparison" ----
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Cells.Item(2, 4).Value = 2
$ws1.Cells.Item(2, 8).Value = 29.05
$ws1.Cells.Item(2, 12).Value = 1.03

# Row 3
$ws1.Cells.Item(3, 4).Value = 3
$ws1.Cells.Item(3, 8).Value = 19
$ws1.Cells.Item(3, 12).Value = 1.17

# Row 4
$ws1.Cells.Item(4, 4).Value = 3
$ws1.Cells.Item(4, 8).Value = 18
$ws1.Cells.Item(4, 12).Value = 0.86

# Row 5
$ws1.Cells.Item(5, 4).Value = 2
$ws1.Cells.Item(5, 8).Value = 25.1
$ws1.Cells.Item(5, 12).Value = 1.08

# Row 6
$ws1.Cells.Item(6, 8).Value = 38.92
$ws1.Cells.Item(6, 12).Value = 0.84

# Row 7
$ws1.Cells.Item(7, 4).Value = 1
$ws1.Cells.Item(7, 8).Value = 37.92
$ws1.Cells.Item(7, 12).Value = 1.14

# Row 8
$ws1.Cells.Item(8, 4).Value = 2
$ws1.Cells.Item(8, 8).Value = 20.87
$ws1.Cells.Item(8, 12).Value = 1.13

# Row 9
$ws1.Cells.Item(9, 8).Value = 13.85
$ws1.Cells.Item(9, 12).Value = 1.16

# Row 10
$ws1.Cells.Item(10, 4).Value = 3
$ws1.Cells.Item(10, 8).Value = 13.68
$ws1.Cells.Item(10, 12).Value = 1.05

# Row 11
$ws1.Cells.Item(11, 8).Value = 20.68
$ws1.Cells.Item(11, 12).Value = 0.85

# Row 12
$ws1.Cells.Item(12, 8).Value = 28.77
$ws1.Cells.Item(12, 12).Value = 0.97

# Row 13
$ws1.Cells.Item(13, 4).Value = 1
$ws1.Cells.Item(13, 8).Value = 21.24
$ws1.Cells.Item(13, 12).Value = 0.93

# Row 14
$ws1.Cells.Item(14, 8).Value = 11.1

# Row 15
$ws1.Cells.Item(15, 4).Value = 3
$ws1.Cells.Item(15, 8).Value = 8.94
$ws1.Cells.Item(15, 12).Value = 0.97

# Row 16
$ws1.Cells.Item(16, 4).Value = 2
$ws1.Cells.Item(16, 8).Value = 10.3
$ws1.Cells.Item(16, 12).Value = 1.04

# Row 17
$ws1.Cells.Item(17, 8).Value = 16.73
$ws1.Cells.Item(17, 12).Value = 0.89

# ---- Sheet 2: "Summary" ----
$ws2 = $wb.Worksheets.Item("Summary")

# Leading apostrophe forces these numeric-looking values to stay text,
# matching the original inline-string cell type.
$ws2.Range("B9").Value = "'37"
$ws2.Range("B10").Value = "'19"
$ws2.Range("B11").Value = "'10"
